$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.UsedRange.EntireColumn.AutoFit() | Out-Null
